$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.831.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.214.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.629"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0934"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.78%  "

$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.548.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.227.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.798.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0959"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.14%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.15%  "

$ws.Range("E28").Value = "  -3.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.81%  "

$ws.Range("E30").Value = "  -1.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0800"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.06%  "

$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.122"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.59%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.90%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0307"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.201"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0997"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("E49").Value = "  -0.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.04%  "
